$d = $word.ActiveDocument

# Original paragraph 1 is a single run "Ghbdtn" followed by the (collapsed)
# "_GoBack" bookmark at the end of the paragraph. The target edit:
#   - changes the existing run's text to "Lehfrb " (note trailing space)
#   - appends a new run "Ghbdtn" right after the bookmark
# i.e. the paragraph text becomes "Lehfrb Ghbdtn" with the bookmark sitting
# between the two runs.

$p1 = $d.Paragraphs(1)
$paraEnd = $p1.Range.End  # position of the paragraph mark

# 1) Append a new run with the original word "Ghbdtn" at the very end of the
#    paragraph (right before the paragraph mark), while the paragraph still
#    only contains "Ghbdtn" - this keeps it a clean, un-"preserve"-flagged
#    run matching the target <w:t>Ghbdtn</w:t>.
$tail = $d.Range($paraEnd - 1, $paraEnd - 1)
$tail.InsertAfter("Ghbdtn")

# 2) Replace only the first (original) occurrence of "Ghbdtn" with
#    "Lehfrb " (trailing space preserved).
$firstWord = $d.Range(0, 6)
$firstWord.Find.Execute("Ghbdtn", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Lehfrb ", 1)

# 3) The engine re-anchors the "_GoBack" bookmark to the end of the
#    paragraph text whenever the paragraph's runs are edited, so restore it
#    to sit right between the two runs (i.e. right after "Lehfrb ").
$bookmarkPos = 7
$d.Bookmarks("_GoBack").Delete()
$d.Bookmarks.Add("_GoBack", $d.Range($bookmarkPos, $bookmarkPos))
